$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows {
    param($ws, $r1, $r2)

    $range1 = $ws.Range("B$r1`:AD$r1")
    $range2 = $ws.Range("B$r2`:AD$r2")

    $vals1 = $range1.Value2
    $vals2 = $range2.Value2

    $range1.Value2 = $vals2
    $range2.Value2 = $vals1
}

Swap-Rows $ws 22 23
Swap-Rows $ws 32 33
Swap-Rows $ws 37 38
Swap-Rows $ws 120 121
Swap-Rows $ws 141 142
Swap-Rows $ws 196 197
Swap-Rows $ws 260 261
Swap-Rows $ws 278 279
Swap-Rows $ws 289 290
Swap-Rows $ws 294 295
Swap-Rows $ws 296 297
